$d = $word.ActiveDocument

# The log table's last row ("January 1st 2022") has three cells:
#   1) Date   2) Awais   3) Zubair
# Awais's cell is currently empty; Zubair's cell already reads
# "- Continued working on flow chart". Both collaborators add a note
# crediting the other.

$t = $d.Tables.Item(1)
$lastRow = $t.Rows.Count

# 1) Awais cell (column 2) is empty -> give it a first run of text.
$cellAwais = $t.Cell($lastRow, 2)
$cellAwais.Range.InsertBefore("- Continued working on flow chart with Zubair")

# 2) Zubair cell (column 3) already has one paragraph -> append a new
#    paragraph with the extra note, keeping the existing paragraph intact.
$cellZubair = $t.Cell($lastRow, 3)
$rZubair = $cellZubair.Range
# Exclude the trailing cell-mark/paragraph-mark character from the range
# so the new paragraph break lands inside the cell, after the existing text.
$rZubair.End = $rZubair.End - 1
$newLine = [char]13
$rZubair.InsertAfter($newLine + "- Continued working on flow chart with Awais")
